# Add 2022-Q3 data
# 1. Insert a new worksheet "2022-Q3" right before "2022-Q2" (i.e. right after "总计"),
#    cloned-in-style from the existing per-quarter fund sheets, populated with the new data.
# 2. Insert a new row at the top of the "总计" summary sheet for 2022-Q3, shifting the
#    existing quarters down by one row, and keep each existing quarter's own figures intact.

$wb = $excel.ActiveWorkbook
$sheets = $wb.Worksheets

# ---------------------------------------------------------------------------
# 1) Create the new "2022-Q3" sheet positioned right before "2022-Q2"
# ---------------------------------------------------------------------------
$q2sheet = $sheets.Item("2022-Q2")
$newSheet = $sheets.Add($q2sheet, $null)
$newSheet.Name = "2022-Q3"

# Header row (B1:H1), matching the style used on the other quarterly fund sheets:
# bold font, thin border all around, centered / top-aligned.
$newSheet.Cells.Item(1,2).Value = "基金代码"
$newSheet.Cells.Item(1,3).Value = "基金名称"
$newSheet.Cells.Item(1,4).Value = "基金规模"
$newSheet.Cells.Item(1,5).Value = "股票总仓位"
$newSheet.Cells.Item(1,6).Value = "仓位占比"
$newSheet.Cells.Item(1,7).Value = "持有市值(亿元)"
$newSheet.Cells.Item(1,8).Value = "仓位排名"

$hdr = $newSheet.Range("B1:H1")
$hdr.Font.Bold = $true
$hdr.Borders.LineStyle = 1
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4160

# Data row (A2:H2)
$newSheet.Cells.Item(2,1).Value = 0
$newSheet.Cells.Item(2,2).Value = "'460010"
$newSheet.Cells.Item(2,3).Value = "华泰柏瑞亚洲领导企业混合（QDII）"
$newSheet.Cells.Item(2,4).Value = "'0.36"
$newSheet.Cells.Item(2,5).Value = "'93.44"
$newSheet.Cells.Item(2,6).Value = "'6.47"
$newSheet.Cells.Item(2,7).Value = "'0.0233"
$newSheet.Cells.Item(2,8).Value = 4

$a2 = $newSheet.Cells.Item(2,1)
$a2.Font.Bold = $true
$a2.Borders.LineStyle = 1
$a2.HorizontalAlignment = -4108
$a2.VerticalAlignment = -4160

# ---------------------------------------------------------------------------
# 2) Insert a new row into "总计" for 2022-Q3 and shift the rest down
# ---------------------------------------------------------------------------
$zj = $sheets.Item("总计")
$zj.Rows.Item(2).Insert()

# the inherited row formatting needs to be cleared on B:D, and A2 needs the
# bold/bordered "index" style used by the rest of column A
$zj.Range("B2:D2").ClearFormats()

$zjA2 = $zj.Cells.Item(2,1)
$zjA2.Font.Bold = $true
$zjA2.Borders.LineStyle = 1
$zjA2.HorizontalAlignment = -4108
$zjA2.VerticalAlignment = -4160

$zj.Cells.Item(2,1).Value = 0
$zj.Cells.Item(2,2).Value = "2022-Q3"
$zj.Cells.Item(2,3).Value = 1
$zj.Cells.Item(2,4).Value = 0.02

# bump the running index in column A for every pre-existing quarter row, now
# shifted one row down (rows 3..8)
for ($r = 3; $r -le 8; $r++) {
    $zj.Cells.Item($r,1).Value = $r - 2
}

Write-Host "done"
